$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the styling of the existing header cells (e.g. H1): bold font,
# thin border on all sides, centered horizontally and top-aligned vertically
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").HorizontalAlignment = -4108
$ws.Range("I1:J1").VerticalAlignment = -4160
$ws.Range("I1:J1").Borders.LineStyle = 1

# Data for columns I (I0) and J (IF), one row per data row (rows 2-65)
$data = @(
    @(4,5),
    @(8,9),
    @(8,9),
    @(9,9),
    @(7,8),
    @(7,8),
    @(8,9),
    @(9,9),
    @(8,9),
    @(8,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(8,8),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(8,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(8,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(7,7),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,10),
    @(9,9),
    @(9,9),
    @(9,9),
    @(8,8),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(6,6),
    @(5,5),
    @(7,7),
    @(5,5)
)

for ($idx = 0; $idx -lt $data.Count; $idx++) {
    $row = $idx + 2
    $pair = $data[$idx]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
